$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @('Enviado', '+55 69 984096742', 'LUCIANO KOVALCHUK - LH 78 SUL, KM 1, LT 29C , GB 11, STR BR 429 SENT/SERING - SAO MIGUEL DO GUAPORE - SAO MIGUEL DO GUAPORE || LUCIANO KOVALCHUK - LH 78 SUL, KM 1, LT 29C , GB 11, STR BR 429 SENT/SERING - SAO MIGUEL DO GUAPORE - SAO MIGUEL DO GUAPORE'),
    @('Enviado', '+55 69 984101481', 'NIRALVA DE OLIVEIRA SILVA - LH 74 SUL, KM 8,5, LT 90, GB 9, STR BR 429 SENT/SERING, DESMEMBRADO - SAO MIGUEL DO GUAPORE - SAO MIGUEL DO GUAPORE'),
    @('Enviado', '+55 69 984396694', 'JESSICA VENANCIO DA SILVA - LH 98 NORTE, KM 5,5, LT 19, GB 16, STR RO 481 - SAO MIGUEL DO GUAPORE - SAO MIGUEL DO GUAPORE'),
    @('Enviado', '+55 69 984421186', 'CLARICE DE FATIMA SOARES NUNES - BR 429, KM 5, LT 33, GB 11, STR BR 429 SENT/SERING - SAO MIGUEL DO GUAPORE - SAO MIGUEL DO GUAPORE'),
    @('Enviado', '+55 69 984432856', 'MANNEZA JORDANIA BERNARDES DE OLIVEIRA - LH 82 SUL, KM 6, LT 33, GB 5, STR RO 481, DESMEMBRADO - SAO MIGUEL DO GUAPORE - SAO MIGUEL DO GUAPORE || MANNEZA JORDANIA BERNARDES DE OLIVEIRA - LH 82 SUL, KM 6, LT 33, GB 5, STR RO 481, DESMEMBRADO - SAO MIGUEL DO GUAPORE - SAO MIGUEL DO GUAPORE'),
    @('Enviado', '+55 69 984735692', 'LUIS DEOTI - BR 429, KM 2, LT 06, GB 11, STR BR 429 SENT/ALV - SAO MIGUEL DO GUAPORE - SAO MIGUEL DO GUAPORE || LUIS DEOTI - BR 429, KM 2, LT 06, GB 11, STR BR 429 SENT/ALV - SAO MIGUEL DO GUAPORE - SAO MIGUEL DO GUAPORE || LUIS DEOTI - BR 429, KM 2, LT 06, GB 11, STR BR 429 SENT/ALV - SAO MIGUEL DO GUAPORE - SAO MIGUEL DO GUAPORE'),
    @('Enviado', '+55 69 984835481', 'SIDNEI BARBOSA SIMAS - LH 25 (MORTA), KM 1,5, LT 05, GB 05, STR BR 429 SENT/SERING - SAO MIGUEL DO GUAPORE - SAO MIGUEL DO GUAPORE'),
    @('Enviado', '+55 69 984956474', 'REMY CARDOSO XAVIER - LH 94 NORTE, KM 3,5, LT 25A, GB 9, STR RO 481 - SAO MIGUEL DO GUAPORE - SAO MIGUEL DO GUAPORE'),
    @('Enviado', '+55 69 992035882', 'JUCELIO DA SILVA PLASTER - LH 90 SUL, KM 5, LT 39, GB 3, STR RO 481 - SAO MIGUEL DO GUAPORE - SAO MIGUEL DO GUAPORE'),
    @('Enviado', '+55 69 992240030', 'MILTON ROSA DE ALMEIDA - BR 429, KM 25, LT 22, STR BR 429 SENT/ALV, DESMENBRADO - SAO MIGUEL DO GUAPORE - SAO MIGUEL DO GUAPORE'),
    @('Enviado', '+55 69 992362662', 'VALDOIR LENZ PEREIRA - BR 429, KM 22, LT 99, GB 01, STR BR 429 SENT/SERING - SAO MIGUEL DO GUAPORE - SAO MIGUEL DO GUAPORE || VALDOIR LENZ PEREIRA - BR 429, KM 22, LT 99, GB 01, STR BR 429 SENT/SERING - SAO MIGUEL DO GUAPORE - SAO MIGUEL DO GUAPORE'),
    @('Enviado', '+55 69 992526904', 'CLERIO MATOS DE ASSIS - LH 78 SUL, KM 3,5, LT 38, GB 11, STR BR 429 SENT/SERING - SAO MIGUEL DO GUAPORE - SAO MIGUEL DO GUAPORE'),
    @('Enviado', '+55 69 992739233', 'DANDARA HORANA DE OLIVEIRA GORZA - LH 86 SUL, KM 11, LT 85, GB 04, STR RO 481, L.D - SAO MIGUEL DO GUAPORE - SAO MIGUEL DO GUAPORE'),
    @('Enviado', '+55 69 992759904', 'ANTONIO MARCOS GOMES - BR 429, KM 9,5, LT 22, GB 10, STR BR 429 SENT/SERING, L.D ACESSO IGREJA CRISTA - SAO MIGUEL DO GUAPORE - SAO MIGUEL DO GUAPORE'),
    @('Enviado', '+55 69 992832220', 'GERALDO GONÇALVES BRANDAO - LH 74 NORTE, KM 4,5, LT 36, GB 14, STR BR 429 SENT/SERING, 2 TRAVESSAO - SAO MIGUEL DO GUAPORE - SAO MIGUEL DO GUAPORE'),
    @('Enviado', '+55 69 992836615', 'OSMAR FERNANDO BORGES - BR 429, KM 23, LT 263, GB 01, STR BR 429 SENT/SERING - SAO MIGUEL DO GUAPORE - SAO MIGUEL DO GUAPORE'),
    @('Enviado', '+55 69 992946955', 'PEDRO PAULO BOA MORTE VAZ - LH 09 SUL, KM 8, LT 51, GB 001, STR BR 429 SENT/SERING - SAO MIGUEL DO GUAPORE - SAO MIGUEL DO GUAPORE'),
    @('Enviado', '+55 69 992972422', 'DALVACY DA PENHA - BR 429, KM 8, LT 46, GB 11, STR BR 429 SENT/SERING - SAO MIGUEL DO GUAPORE - SAO MIGUEL DO GUAPORE'),
    @('Enviado', '+55 69 993039887', 'GILMAR ANDERSON DA SILVA - LH 98 SUL, KM 12, LT 12, GB 10, STR RO 481, ROD RO 481 - SAO MIGUEL DO GUAPORE - SANTANA DO GUAPORE'),
    @('Enviado', '+55 69 993139237', 'ATAIDE PEREIRA DE SOUZA - LH 82 SUL, KM 6, LT 33, GB 5, STR RO 481, ENFRENTE O TANQUE DE LEITE - SAO MIGUEL DO GUAPORE - SAO MIGUEL DO GUAPORE'),
    @('Enviado', '+55 69 993185430', 'EILANE GONÇALVES DE ALMEIDA - LH 86 SUL, KM 10, LT 86, GB 3, STR RO 481 - SAO MIGUEL DO GUAPORE - SAO MIGUEL DO GUAPORE'),
    @('Enviado', '+55 69 993207066', 'VALDECI TIMOTEO - LH 90 SUL, KM 7, LT 60-A3, GB 2, STR RO 481 - SAO MIGUEL DO GUAPORE - SAO MIGUEL DO GUAPORE'),
    @('Enviado', '+55 69 993266509', 'ORLANDINO ALVES DE OLIVEIRA - LH 94 SUL, KM 9, LT 37, GB 1, STR RO 481 - SAO MIGUEL DO GUAPORE - SAO MIGUEL DO GUAPORE'),
    @('Enviado', '+55 69 993277937', 'DANDARA HORANA DE OLIVEIRA GORZA - LH 86 SUL, KM 11, LT 85, GB 04, STR RO 481, L.D - SAO MIGUEL DO GUAPORE - SAO MIGUEL DO GUAPORE'),
    @('Número inválido', '+55 69 993317270', 'SIRLENE LIZZI CAVALCANTE - LH 86 SUL, KM 19, LT 149, 149-A, GB 004, STR RO 481 - SAO MIGUEL DO GUAPORE - SAO MIGUEL DO GUAPORE'),
    @('Número inválido', '+55 69 993449821', 'EILANE GONÇALVES DE ALMEIDA - LH 86 SUL, KM 10, LT 86, GB 3, STR RO 481 - SAO MIGUEL DO GUAPORE - SAO MIGUEL DO GUAPORE'),
    @('Enviado', '+55 69 993617963', 'ISAIAS STRELOW - LH 98 NORTE, KM 10, LT 34-A, GB 16, STR RO 481, DESMENBRADO - SAO MIGUEL DO GUAPORE - SAO MIGUEL DO GUAPORE'),
    @('Enviado', '+55 69 993780921', 'JUCELIO DA SILVA PLASTER - LH 90 SUL, KM 5, LT 39, GB 3, STR RO 481 - SAO MIGUEL DO GUAPORE - SAO MIGUEL DO GUAPORE'),
    @('Enviado', '+55 69 993860390', 'RAIMUNDO DA CONCEIÇAO SA - BR 429, KM 7,5, LT 20, GB 14, STR BR 429 SENT/SERING, SETOR SAO MIGUEL - SAO MIGUEL DO GUAPORE - SAO MIGUEL DO GUAPORE'),
    @('Enviado', '+55 69 996096955', 'FELIX SOUZA ESPIRITO SANTO - LH 78 SUL, KM 7, LT 60, GB 05, STR BR 429 SENT/SERING - SAO MIGUEL DO GUAPORE - SAO MIGUEL DO GUAPORE || FELIX SOUZA ESPIRITO SANTO - LH 78 SUL, KM 7, LT 60, GB 05, STR BR 429 SENT/SERING - SAO MIGUEL DO GUAPORE - SAO MIGUEL DO GUAPORE'),
    @('Enviado', '+55 69 999019031', 'VALDECI TIMOTEO - LH 90 SUL, KM 7, LT 60-A3, GB 2, STR RO 481 - SAO MIGUEL DO GUAPORE - SAO MIGUEL DO GUAPORE'),
    @('Enviado', '+55 69 999126066', 'FRANCISCO ALVES DA SILVA - LH 90 SUL, KM 16, LT 130-B, GB 02, STR RO 481 - SAO MIGUEL DO GUAPORE - SAO MIGUEL DO GUAPORE'),
    @('Enviado', '+55 69 999232663', 'BRUNO FARIAS CRAUSE - LH 94 NORTE, KM 4, LT 29, GB 09, STR RO 481 - SAO MIGUEL DO GUAPORE - SAO MIGUEL DO GUAPORE'),
    @('Enviado', '+55 69 999268048', 'GEISILANDE BRAZ DUARTE - LH 90 SUL, KM 2,5, LT 22, GB 02, STR RO 481, DESMEMBRADO - SAO MIGUEL DO GUAPORE - SAO MIGUEL DO GUAPORE'),
    @('Enviado', '+55 69 999419677', 'WALLAFE DO NASCIMENTO CAETANO - LH 86 SUL, KM 5, LT 41, GB 04, STR RO 481, DESMEMBRADO - SAO MIGUEL DO GUAPORE - SAO MIGUEL DO GUAPORE'),
    @('Enviado', '+55 69 999455742', 'APARECIDO ANTONIO DE SOUZA - LH 94 SUL, KM 4, LT 029, GB 002, STR RO 481 - SAO MIGUEL DO GUAPORE - SAO MIGUEL DO GUAPORE || APARECIDO ANTONIO DE SOUZA - LH 94 SUL, KM 4, LT 029, GB 002, STR RO 481 - SAO MIGUEL DO GUAPORE - SAO MIGUEL DO GUAPORE'),
    @('Enviado', '+55 69 999612659', 'GEISILANDE BRAZ DUARTE - LH 90 SUL, KM 2,5, LT 22, GB 02, STR RO 481, DESMEMBRADO - SAO MIGUEL DO GUAPORE - SAO MIGUEL DO GUAPORE'),
    @('Enviado', '+55 69 999621815', 'WILSON ROSA DE ALMEIDA - BR 429, KM 25, LT 22-F, GB 05, STR BR 429 SENT/ALV, DESMEMBRADO - SAO MIGUEL DO GUAPORE - SAO MIGUEL DO GUAPORE'),
    @('Enviado', '+55 69 999804992', 'VALDEIR TIMOTEO - LH 90 SUL, KM 7, LT 60, GB 2, STR RO 481 - SAO MIGUEL DO GUAPORE - SAO MIGUEL DO GUAPORE'),
    @('Enviado', '+55 69 99996857', 'GILMAR BELIZARIO DA SILVA - LH 78 SUL, KM 5, LT 35 B, GB 06, STR BR 429 SENT/SERING - SAO MIGUEL DO GUAPORE - SAO MIGUEL DO GUAPORE')
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $r++
}